# TFs stable, RPT development
# Swap the transfer-function text in C21/C22 on the "Parameters" sheet:
#   C21: tfIfdl(s,[0 2:3],%s)  ->  tfIfdl(s,[2:3],%s)
#   C22: tfPhie(s,[0:3],%s)    ->  tfPhie(s,[1:3],%s)
# and move the active selection to G23.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")

# Leading apostrophe keeps these as forced-text entries (matches the
# existing quote-prefixed cell formatting) without becoming part of the
# stored string.
$ws.Range("C21").Value = "tfIfdl(s,[2:3],%s)"
$ws.Range("C22").Value = "'tfPhie(s,[1:3],%s)"

$ws.Range("G23").Select()
